$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Configuration")
$ws.Range("C8").Value = "DICOM:StudyDate"
$ws.Range("C8").Select()
